$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Paragraph 6: "    Encourage these neighbors on the rendered depth of NeRF to be continuous. ..." ---
# Color the substring "on the rendered depth of NeRF " (30 chars right after the
# 30-char lead-in "    Encourage these neighbors ") in blue (0432FF).
$para6 = $tr.Paragraphs(6)
$highlight6 = $tr.Characters($para6.Start + 30, 30)
$highlight6.Font.Color.RGB = 16724484

# --- Paragraph 7: "           we compute the loss2= loss2+ max(|depth_A-depth_A_i|-m', 0)" ---
# Turn "depth_A-depth_A_i" into "depth'_A-depth'_A_i" (insert a right single
# quote after each "depth"), keeping the run split before the trailing "A_i".
$para7 = $tr.Paragraphs(7)
$run1 = $tr.Characters($para7.Start, 45)
$run1.Text = "           we compute the loss2= loss2+ max(|depth" + [char]0x2019 + "_A-depth" + [char]0x2019 + "_"
$run2 = $tr.Characters($para7.Start + 61, 17)
$run2.Text = "A_i"
